# Add a new row 12 ("2021年") to Sheet1, mirroring the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous data row (row 11) down into row 12
# so the new row picks up the same cell style (bold/border/center for A12).
$ws.Range("A11:U11").Copy()
$ws.Range("A12:U12").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Cells.Item(12, 1).Value = "2021年"
$ws.Cells.Item(12, 2).Value = 1092162.7
$ws.Cells.Item(12, 4).Value = 1046412.6
$ws.Cells.Item(12, 5).Value = 4414577.9
$ws.Cells.Item(12, 6).Value = 1816552.3
$ws.Cells.Item(12, 7).Value = 4850510
$ws.Cells.Item(12, 8).Value = 9424368.4
$ws.Cells.Item(12, 9).Value = 1331366.4
$ws.Cells.Item(12, 10).Value = 421734.9
$ws.Cells.Item(12, 11).Value = 2715179.8
$ws.Cells.Item(12, 13).Value = 36345892.3
$ws.Cells.Item(12, 14).Value = 10323259.7
$ws.Cells.Item(12, 15).Value = 3163875.9
$ws.Cells.Item(12, 16).Value = 506005.6
$ws.Cells.Item(12, 17).Value = 1008550.2
$ws.Cells.Item(12, 20).Value = 238006
$ws.Cells.Item(12, 21).Value = 56845724.2
